$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 3-32 hold District data in column G.
# Row 17 and 26 previously had an empty F cell that spilled text into G;
# after the correction, F is cleared (no cell) and G holds the plain name.
$special = @{
    17 = "Vijayapura"
    26 = "Vijayapura (Bijapur)"
}

for ($row = 3; $row -le 32; $row++) {
    if ($special.ContainsKey($row)) {
        $ws.Cells.Item($row, 7).Value = $special[$row]
    } else {
        $ws.Cells.Item($row, 7).Value = "Vijayapura (Bijapur)"
    }
}

# Clear the stray empty inline-string cells in column F for rows 17 and 26
$ws.Cells.Item(17, 6).Clear()
$ws.Cells.Item(26, 6).Clear()
